# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-22) were sorted descending (1903 -> 1809);
# this update re-sorts them ascending (1809 -> 1903), keeping each
# period's "Valor Mora" value attached to the same period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order of (period, value) pairs, row by row (16..22).
$periods = @("1809", "1810", "1811", "1812", "1901", "1902", "1903")
$values  = @(18749, 23437, 23437, 23437, 31249, 31249, 31249)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $values[$i]
}
